# Auto update Excel log
# Appends rows 46-58 (2026-01-28 sensor readings) to the PIR, Humidity,
# and Temperature sheets. Column A (Date) and the Humidity sheet's
# column E (percentage-looking text) are forced to Text format first so
# Excel's autodetection does not silently convert them to date/number
# serials - the source log stores every value as literal text.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A46:A58").NumberFormat = "@"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("B46").Value = "17:05:43"
$ws.Range("C46").Value = "17:00"
$ws.Range("D46").Value = "Bathroom"
$ws.Range("E46").Value = "No Motion"
$ws.Range("F46").Value = "Inactive"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("B47").Value = "17:05:44"
$ws.Range("C47").Value = "17:00"
$ws.Range("D47").Value = "Bathroom"
$ws.Range("E47").Value = "No Motion"
$ws.Range("F47").Value = "Inactive"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("B48").Value = "17:05:48"
$ws.Range("C48").Value = "17:00"
$ws.Range("D48").Value = "Bathroom"
$ws.Range("E48").Value = "No Motion"
$ws.Range("F48").Value = "Inactive"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("B49").Value = "17:05:53"
$ws.Range("C49").Value = "17:00"
$ws.Range("D49").Value = "Bathroom"
$ws.Range("E49").Value = "No Motion"
$ws.Range("F49").Value = "Inactive"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("B50").Value = "17:05:58"
$ws.Range("C50").Value = "17:00"
$ws.Range("D50").Value = "Bathroom"
$ws.Range("E50").Value = "No Motion"
$ws.Range("F50").Value = "Inactive"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("B51").Value = "17:06:03"
$ws.Range("C51").Value = "17:00"
$ws.Range("D51").Value = "Bathroom"
$ws.Range("E51").Value = "No Motion"
$ws.Range("F51").Value = "Inactive"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("B52").Value = "17:06:08"
$ws.Range("C52").Value = "17:00"
$ws.Range("D52").Value = "Bathroom"
$ws.Range("E52").Value = "No Motion"
$ws.Range("F52").Value = "Inactive"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("B53").Value = "17:06:13"
$ws.Range("C53").Value = "17:00"
$ws.Range("D53").Value = "Bathroom"
$ws.Range("E53").Value = "No Motion"
$ws.Range("F53").Value = "Inactive"
$ws.Range("A54").Value = "2026-01-28"
$ws.Range("B54").Value = "17:06:18"
$ws.Range("C54").Value = "17:00"
$ws.Range("D54").Value = "Bathroom"
$ws.Range("E54").Value = "No Motion"
$ws.Range("F54").Value = "Inactive"
$ws.Range("A55").Value = "2026-01-28"
$ws.Range("B55").Value = "17:06:23"
$ws.Range("C55").Value = "17:00"
$ws.Range("D55").Value = "Bathroom"
$ws.Range("E55").Value = "No Motion"
$ws.Range("F55").Value = "Inactive"
$ws.Range("A56").Value = "2026-01-28"
$ws.Range("B56").Value = "17:06:28"
$ws.Range("C56").Value = "17:00"
$ws.Range("D56").Value = "Bathroom"
$ws.Range("E56").Value = "No Motion"
$ws.Range("F56").Value = "Inactive"
$ws.Range("A57").Value = "2026-01-28"
$ws.Range("B57").Value = "17:06:33"
$ws.Range("C57").Value = "17:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "No Motion"
$ws.Range("F57").Value = "Inactive"
$ws.Range("A58").Value = "2026-01-28"
$ws.Range("B58").Value = "17:06:38"
$ws.Range("C58").Value = "17:00"
$ws.Range("D58").Value = "Bathroom"
$ws.Range("E58").Value = "No Motion"
$ws.Range("F58").Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A46:A58").NumberFormat = "@"
$ws.Range("E46:E58").NumberFormat = "@"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("B46").Value = "17:05:43"
$ws.Range("C46").Value = "17:00"
$ws.Range("D46").Value = "Bathroom"
$ws.Range("E46").Value = "87.6%"
$ws.Range("F46").Value = "Active"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("B47").Value = "17:05:44"
$ws.Range("C47").Value = "17:00"
$ws.Range("D47").Value = "Bathroom"
$ws.Range("E47").Value = "86.7%"
$ws.Range("F47").Value = "Active"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("B48").Value = "17:05:47"
$ws.Range("C48").Value = "17:00"
$ws.Range("D48").Value = "Bathroom"
$ws.Range("E48").Value = "87.6%"
$ws.Range("F48").Value = "Active"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("B49").Value = "17:05:51"
$ws.Range("C49").Value = "17:00"
$ws.Range("D49").Value = "Bathroom"
$ws.Range("E49").Value = "86.7%"
$ws.Range("F49").Value = "Active"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("B50").Value = "17:05:55"
$ws.Range("C50").Value = "17:00"
$ws.Range("D50").Value = "Bathroom"
$ws.Range("E50").Value = "87.6%"
$ws.Range("F50").Value = "Active"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("B51").Value = "17:06:04"
$ws.Range("C51").Value = "17:00"
$ws.Range("D51").Value = "Bathroom"
$ws.Range("E51").Value = "86.7%"
$ws.Range("F51").Value = "Active"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("B52").Value = "17:06:16"
$ws.Range("C52").Value = "17:00"
$ws.Range("D52").Value = "Bathroom"
$ws.Range("E52").Value = "87.7%"
$ws.Range("F52").Value = "Active"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("B53").Value = "17:06:20"
$ws.Range("C53").Value = "17:00"
$ws.Range("D53").Value = "Bathroom"
$ws.Range("E53").Value = "87.6%"
$ws.Range("F53").Value = "Active"
$ws.Range("A54").Value = "2026-01-28"
$ws.Range("B54").Value = "17:06:24"
$ws.Range("C54").Value = "17:00"
$ws.Range("D54").Value = "Bathroom"
$ws.Range("E54").Value = "86.7%"
$ws.Range("F54").Value = "Active"
$ws.Range("A55").Value = "2026-01-28"
$ws.Range("B55").Value = "17:06:28"
$ws.Range("C55").Value = "17:00"
$ws.Range("D55").Value = "Bathroom"
$ws.Range("E55").Value = "87.7%"
$ws.Range("F55").Value = "Active"
$ws.Range("A56").Value = "2026-01-28"
$ws.Range("B56").Value = "17:06:32"
$ws.Range("C56").Value = "17:00"
$ws.Range("D56").Value = "Bathroom"
$ws.Range("E56").Value = "86.7%"
$ws.Range("F56").Value = "Active"
$ws.Range("A57").Value = "2026-01-28"
$ws.Range("B57").Value = "17:06:36"
$ws.Range("C57").Value = "17:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "87.6%"
$ws.Range("F57").Value = "Active"
$ws.Range("A58").Value = "2026-01-28"
$ws.Range("B58").Value = "17:06:40"
$ws.Range("C58").Value = "17:00"
$ws.Range("D58").Value = "Bathroom"
$ws.Range("E58").Value = "87.6%"
$ws.Range("F58").Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A46:A58").NumberFormat = "@"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("B46").Value = "17:05:44"
$ws.Range("C46").Value = "17:00"
$ws.Range("D46").Value = "Bathroom"
$ws.Range("E46").Value = "22.9C"
$ws.Range("F46").Value = "Active"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("B47").Value = "17:05:44"
$ws.Range("C47").Value = "17:00"
$ws.Range("D47").Value = "Bathroom"
$ws.Range("E47").Value = "22.9C"
$ws.Range("F47").Value = "Active"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("B48").Value = "17:05:48"
$ws.Range("C48").Value = "17:00"
$ws.Range("D48").Value = "Bathroom"
$ws.Range("E48").Value = "22.9C"
$ws.Range("F48").Value = "Active"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("B49").Value = "17:05:52"
$ws.Range("C49").Value = "17:00"
$ws.Range("D49").Value = "Bathroom"
$ws.Range("E49").Value = "22.9C"
$ws.Range("F49").Value = "Active"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("B50").Value = "17:05:56"
$ws.Range("C50").Value = "17:00"
$ws.Range("D50").Value = "Bathroom"
$ws.Range("E50").Value = "22.8C"
$ws.Range("F50").Value = "Active"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("B51").Value = "17:06:04"
$ws.Range("C51").Value = "17:00"
$ws.Range("D51").Value = "Bathroom"
$ws.Range("E51").Value = "22.9C"
$ws.Range("F51").Value = "Active"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("B52").Value = "17:06:16"
$ws.Range("C52").Value = "17:00"
$ws.Range("D52").Value = "Bathroom"
$ws.Range("E52").Value = "22.9C"
$ws.Range("F52").Value = "Active"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("B53").Value = "17:06:20"
$ws.Range("C53").Value = "17:00"
$ws.Range("D53").Value = "Bathroom"
$ws.Range("E53").Value = "22.9C"
$ws.Range("F53").Value = "Active"
$ws.Range("A54").Value = "2026-01-28"
$ws.Range("B54").Value = "17:06:24"
$ws.Range("C54").Value = "17:00"
$ws.Range("D54").Value = "Bathroom"
$ws.Range("E54").Value = "22.9C"
$ws.Range("F54").Value = "Active"
$ws.Range("A55").Value = "2026-01-28"
$ws.Range("B55").Value = "17:06:28"
$ws.Range("C55").Value = "17:00"
$ws.Range("D55").Value = "Bathroom"
$ws.Range("E55").Value = "22.9C"
$ws.Range("F55").Value = "Active"
$ws.Range("A56").Value = "2026-01-28"
$ws.Range("B56").Value = "17:06:32"
$ws.Range("C56").Value = "17:00"
$ws.Range("D56").Value = "Bathroom"
$ws.Range("E56").Value = "22.9C"
$ws.Range("F56").Value = "Active"
$ws.Range("A57").Value = "2026-01-28"
$ws.Range("B57").Value = "17:06:36"
$ws.Range("C57").Value = "17:00"
$ws.Range("D57").Value = "Bathroom"
$ws.Range("E57").Value = "22.8C"
$ws.Range("F57").Value = "Active"
$ws.Range("A58").Value = "2026-01-28"
$ws.Range("B58").Value = "17:06:40"
$ws.Range("C58").Value = "17:00"
$ws.Range("D58").Value = "Bathroom"
$ws.Range("E58").Value = "22.9C"
$ws.Range("F58").Value = "Active"

